$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add column D values (word / exact-phrase match counts) for rows 2-53 ---
$dValues = @{
    2  = 6
    3  = 3
    4  = 11
    5  = 13
    6  = 9
    7  = 36
    8  = 27
    9  = 4
    10 = 18
    11 = 8
    12 = 17
    13 = 15
    14 = 12
    15 = 23
    16 = 29
    17 = 21
    18 = 15
    19 = 55
    20 = 8
    21 = 11
    22 = 14
    23 = 11
    24 = 101
    25 = 24
    26 = 10
    27 = 53
    28 = 12
    29 = 16
    30 = 60
    31 = 27
    32 = 19
    33 = 81
    34 = 38
    35 = 10
    36 = 31
    37 = 16
    38 = 24
    39 = 12
    40 = 15
    41 = 16
    42 = 40
    43 = 14
    44 = 9
    45 = 15
    46 = 11
    47 = 74
    48 = 30
    49 = 15
    50 = 8
    51 = 42
    52 = 32
    53 = 13
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# --- Widen column B so the full names/phrases are visible ---
$ws.Columns("B").ColumnWidth = 50.28515625

# --- Rows whose wrapped-text height needs to shrink now that column B is wider ---
$ws.Rows("56:57").EntireRow.AutoFit()
$ws.Rows("58:59").RowHeight = 30
$ws.Rows("60:61").EntireRow.AutoFit()
$ws.Rows("63:65").EntireRow.AutoFit()
$ws.Rows("67:69").EntireRow.AutoFit()

# --- Update the view: scrolled down a bit, with B16 now selected ---
$ws.Application.Goto($ws.Range("A7"), $false)
$ws.Range("B16").Select()
